$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row 8: financial period labels shift forward by one year ---
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# --- Header row 9: publish dates shift ---
$ws.Range("D9").Value = "1399-04-04 (8)"
$ws.Range("E9").Value = "1400-04-05 (11)"
$ws.Range("F9").Value = "1401-04-19 (12)"
$ws.Range("G9").Value = "1402-02-30 (9)"
$ws.Range("H9").Value = "1402-02-30"

# --- Data rows: reset financial figures (new read_price algorithm pending fresh data) ---
$ws.Range("D11:H11").Value = 0
$ws.Range("D12:H12").Value = 0
$ws.Range("D13:H13").Value = 0
$ws.Range("D14:H14").Value = 0

$ws.Range("E15").Value = "-"
$ws.Range("F15").Value = "-"
$ws.Range("G15").Value = "-"
$ws.Range("H15").Value = "-"

$ws.Range("D16:H16").Value = 0
$ws.Range("D17:H17").Value = 0
$ws.Range("D18:H18").Value = 0
$ws.Range("D19:H19").Value = 0
$ws.Range("D20:H20").Value = 0

$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = "-"
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 0

$ws.Range("D22:H22").Value = 0

$ws.Range("D23").Value = "-"
$ws.Range("E23").Value = "-"
$ws.Range("F23").Value = "-"
$ws.Range("G23").Value = "-"
$ws.Range("H23").Value = "-"

$ws.Range("D24:H24").Value = 0

$ws.Range("D25").Value = "-"
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0

$ws.Range("D26:H26").Value = 0
$ws.Range("D27:H27").Value = 0
